$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = "γ谷氨酰转肽酶"

# Row 7
$ws.Range("C7").Value = "μmol/L"

# Row 9
$ws.Range("C9").Value = "μmol/L"

# Row 10
$ws.Range("C10").Value = "μmol/L"

# Row 11
$ws.Range("A11").Value = "白蛋白"

# Row 13
$ws.Range("A13").Value = "白蛋白"

# Row 14
$ws.Range("A14").Value = "球蛋白"

# Row 15
$ws.Range("A15").Value = "白球比"

# Row 16
$ws.Range("A16").Value = "尿素"

# Row 17
$ws.Range("A17").Value = "肌酐"
$ws.Range("C17").Value = "μmol/L"

# Row 18
$ws.Range("A18").Value = "尿酸"
$ws.Range("C18").Value = "μmol/L"

# Row 22
$ws.Range("B22").Value = "'2.3"

# Row 23
$ws.Range("C23").Value = "mmol/L"

# Row 24
$ws.Range("A24").Value = "镁"
$ws.Range("C24").Value = "mmol/L"

# Row 25
$ws.Range("A25").Value = "总二氧化碳"
$ws.Range("C25").Value = "mmol/L"

# Row 26
$ws.Range("A26").Value = "总胆固醇"

# Row 27
$ws.Range("A27").Value = "甘油三酯"
$ws.Range("C27").Value = "mmol/L"

# Row 28
$ws.Range("A28").Value = "高密度脂蛋白胆固醇"
$ws.Range("B28").Value = "'1.0"
$ws.Range("C28").Value = "mmol/L"

# Row 29
$ws.Range("A29").Value = "低密度脂蛋白胆固醇"

# Row 30
$ws.Range("A30").Value = "载脂蛋白-B"
$ws.Range("C30").Value = "g/L"
$ws.Range("D30").Value = "1.00-175"

# Row 31
$ws.Range("A31").Value = "载脂蛋白-B"
$ws.Range("D31").Value = "0.60-1.10"

# Row 32
$ws.Range("A32").Value = "白蛋白"
$ws.Range("C32").Value = "mg/L"
$ws.Range("D32").Value = "0-300"

# Row 33
$ws.Range("A33").Value = "血糖"
$ws.Range("C33").Value = "mmol/L"
$ws.Range("D33").Value = "3.89-6.11"
